$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '22.490.17'
$ws.Range("E2").Value = '  +0.09%  '

$ws.Range("D3").Value = '1.572.62'
$ws.Range("E3").Value = '  +0.19%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.004'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.09%  '

$ws.Range("E5").Value = '  +0.10%  '

$ws.Range("E6").Value = '  -0.68%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3724'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.71%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '48.33'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -2.97%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3338'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.06%  '

$ws.Range("E10").Value = '  -1.24%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07492'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.50%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.004'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.13%  '

$ws.Range("E13").Value = '  -0.96%  '

$ws.Range("E14").Value = '  -0.68%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.930'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.50%  '

$ws.Range("D16").Value = '1.572.62'
$ws.Range("E16").Value = '  +0.03%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001118'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.21%  '

$ws.Range("E18").Value = '  -2.17%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06779'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.13%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.003'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.08%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.393'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.65%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '16.49'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.48%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '12.12'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.07%  '

$ws.Range("D24").Value = '22.474.19'
$ws.Range("E24").Value = '  +0.03%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.393'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.58%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.577'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.21%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '152.72'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.48%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '19.74'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.36%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.004'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.35%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '124.17'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.61%  '

$ws.Range("D31").Value = '1.753.24'
$ws.Range("E31").Value = '  +0.18%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.054'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.09%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.170'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.24%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.016'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.14%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '9.705'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.96%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.08320'
$ws.Range("D36").Style = "Normal"

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02461'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.80%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2280'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.87%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.06394'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.63%  '

$ws.Range("B40").Value = 'InternetComputer(DFINITY)'
$ws.Range("C40").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.391'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.88%  '

$ws.Range("B41").Value = 'TrustWalletToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.295'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -4.27%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.6324'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.69%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '11.32'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.61%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.003'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.12%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.89'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.45%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.6154'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +5.13%  '

$ws.Range("E47").Value = '  -0.56%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.062'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.43%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '125.50'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.74%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.215'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.77%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.07269'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.47%  '
